## act tablas web jul25
## Adds the 2023 and 2022 data points to the "Data" sheet (shifting the
## existing Fecha/Valor series down by two rows) and records the July 2025
## update note on the "Metadata" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Data" sheet: prepend 2023 / 2022 to the Fecha/Valor time series.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Fecha (years, newest first) and the matching Valor figures, row 2..20.
$years = @("2023","2022","2021","2020","2019","2018","2017","2016","2015","2014","2013","2012","2011","2010","2009","2008","2007","2006","2005")
$vals  = @(14.5,13.9,14.4,13.6,14,13.9,13.7,13.6,13.6,13.9,13.8,14.1,13.9,13.5,13.9,13.4,11.8,10.7,10.5)

# Force column A to be entered as text (so "2023" etc. stay strings, same
# as the existing year labels) instead of being auto-coerced to numbers.
$yearsRange = $ws.Range("A2:A20")
$yearsRange.NumberFormat = "@"

for ($i = 0; $i -lt $years.Length; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = $years[$i]
    $ws.Range("B$r").Value = $vals[$i]
}

# Restore the default (Normal) style now that the text is safely in place.
$yearsRange.Style = "Normal"

# ---------------------------------------------------------------------
# 2) "Metadata" sheet: insert the "actualizacion" / "Julio 2025" row just
#    above the "cita" row, and normalise the leading blank row.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Find the "cita" row so the new row lands in the right place even if the
# sheet layout ever shifts.
$citaRow = 0
for ($r = 1; $r -le $meta.UsedRange.Rows.Count; $r++) {
    if ($meta.Range("A$r").Value() -eq "cita") {
        $citaRow = $r
        break
    }
}

$meta.Rows.Item($citaRow).Insert()
$meta.Range("A$citaRow").Value = "actualizacion"
$meta.Range("B$citaRow").Value = "Julio 2025"

# The first row's A cell was a blank placeholder; align it with the blank
# used elsewhere in the sheet (a single space).
$meta.Range("A1").Value = " "
